# Add new "Flanger" effect parameters to the Generator Scripts Parameters workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Vibrato's "LFO Waveform" enum count changes from 4 to 5 (row 56, column E)
$ws.Range("E56").Value = 5

# New data rows 61-68 describing the "Flanger" effect, modeled on the existing
# Phaser / Vibrato effect rows.
$flangerRows = @(
    @{ Row = 61; B = "Min Delay";         C = "minDelay";   D = "Min Delay";         F = "ms"; G = 1;    H = 20;  I = 2.5; J = 0.5;  K = "minDelayMs" },
    @{ Row = 62; B = "Sweep Width";       C = "sweepWidth"; D = "Sweep Width";       F = "ms"; G = 1;    H = 20;  I = 10;  J = 0.5;  K = "sweepWidthMs" },
    @{ Row = 63; B = "Depth";             C = "depth";      D = "Depth";             F = "%";  G = 0;    H = 100; I = 100; J = 1;    K = "depth" },
    @{ Row = 64; B = "Feedback";          C = "feedback";   D = "Feedback";          F = "%";  G = 0;    H = 50;  I = 0;   J = 1;    K = "feedback" },
    @{ Row = 65; B = "LFO Frequency";     C = "lfoFreq";    D = "LFO Freq";          F = "Hz"; G = 0.05; H = 2;   I = 0.2; J = 0.01; K = "lfoFreqHz" }
)

foreach ($r in $flangerRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "Flanger"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
}

# Rows that only have an enum count (E) instead of min/max/step range.
$flangerEnumRows = @(
    @{ Row = 66; B = "Stereo";             C = "stereo";      D = "Stereo";             E = 2; I = 0; K = "stereoMode" },
    @{ Row = 67; B = "LFO Waveform";       C = "lfoWaveform"; D = "LFO";                E = 4; I = 0; K = "lfoWaveform" },
    @{ Row = 68; B = "Interpolation Type"; C = "interpType";  D = "Interpolation Type"; E = 3; I = 1; K = "interpolationType" }
)

foreach ($r in $flangerEnumRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "Flanger"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 11).Value = $r.K
}

# Match the "A" column style (bold/shaded header style) used for every effect
# name cell, copied from the Vibrato block immediately above.
$ws.Range("A56").Copy()
$ws.Range("A61:A68").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Match the numeric cell style (centered) used in columns E and G:J for the
# other parameter rows, copied from the Vibrato block immediately above.
$ws.Range("E56").Copy()
$ws.Range("E66:E68").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F57:J57").Copy()
$ws.Range("F61:J61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F58:J58").Copy()
$ws.Range("F62:J62").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F39:J39").Copy()
$ws.Range("F63:J63").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F40:J40").Copy()
$ws.Range("F64:J64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F41:J41").Copy()
$ws.Range("F65:J65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the view to reflect where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("D63").Select()
